# Commit: "update- this is json flow commit"
#
# The edit collapses the "Concise Summary" section entirely (it and all of
# its child bullet content is removed) and also strips everything that used
# to live under the "Fulsome Summary" heading, so that only the bare
# "Fulsome Summary" heading paragraph remains - directly after the document
# title and the blank paragraph that follows it. The "Fulsome Summary"
# paragraph itself keeps its original paragraph/run formatting (bold,
# underlined, sz 22, spacing-after 80) untouched.

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

$concisePara = Find-ParagraphByText $d "Concise Summary"
$fulsomePara = Find-ParagraphByText $d "Fulsome Summary"

# 1) Delete the "Concise Summary" heading paragraph plus everything in its
#    section, up to (but not including) the "Fulsome Summary" paragraph.
$sectionRange = $d.Range($concisePara.Range.Start, $fulsomePara.Range.Start)
$sectionRange.Delete()

# 2) Re-locate "Fulsome Summary" (its index shifted after the deletion
#    above) and delete everything from the end of that paragraph through
#    the end of the document body - i.e. the whole "Fulsome Summary"
#    section content.
$fulsomePara = Find-ParagraphByText $d "Fulsome Summary"
$tailRange = $d.Range($fulsomePara.Range.End, $d.Content.End)
$tailRange.Delete()
